# Apply the "copied some cells from template xlsx to invoice xlsx" edit:
#  1. Clear the placeholder name text out of A2 (keeps its formatting/style).
#  2. Remove the blank spacer row 32 (rows below shift up by one).
#  3. Rename the sheet from "Template" to "Factuur".
#  4. Re-point the print area to the renamed sheet / new (smaller) used range.
#  5. Select I37 to match the saved selection/view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. A2 held the literal name "Denise Boezaard" - clear it, leave style s="1".
$ws.Range("A2").ClearContents()

# 2. Delete the empty row 32; rows 33-38 shift up to become 32-37.
$ws.Rows(32).Delete()

# 3. Rename worksheet.
$ws.Name = "Factuur"

# 4. Update the print area (defined name follows the sheet's new name/range).
$ws.PageSetup.PrintArea = '$A$1:$H$36'

# 5. Match the final saved selection / scroll position.
$ws.Range("I37").Select() | Out-Null
